$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new scores must land as literal text (shared strings), matching the
# rest of the table (which already stores its numbers as text). Writing a
# numeric-looking string straight into Range.Value gets auto-coerced into a
# real number by Excel, so instead render it through a TEXT() formula in a
# scratch cell, then Copy/PasteSpecial-Values it into place - this carries
# the text result over without stamping a new (unwanted) number-format style
# on the destination cell.
function Set-TextValue($cellAddr, $text) {
    $scratch = $ws.Range("Z100")
    $scratch.Formula = "=TEXT(" + $text + ",""0.00"")"
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
    $scratch.ClearContents()
}

Set-TextValue "B3" "6787.92"
Set-TextValue "D5" "6788.92"

Set-TextValue "E2" "6646.20"
Set-TextValue "E3" "6806.28"
Set-TextValue "E4" "6812.20"
Set-TextValue "E5" "6803.92"

$excel.CutCopyMode = 0

# Update the selection to match the new active cell recorded in the file.
$ws.Range("E5").Select()
